# Update TPM-derived NATMI metrics for the Oxt-Oxtr ligand-receptor sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.4557103333333333
$ws.Range("H2").Value = 1.367131
$ws.Range("I2").Value = 0.1996401272959883
$ws.Range("J2").Value = 0.1996401272959883
$ws.Range("M2").Value = 0.029575
$ws.Range("N2").Value = 0.088725
$ws.Range("Q2").Value = 0.01347763310833333
$ws.Range("R2").Value = 0.121298697975
$ws.Range("S2").Value = 0.1996401272959883
$ws.Range("T2").Value = 0.1996401272959883

# Row 3
$ws.Range("I3").Value = 0.2962807848215612
$ws.Range("J3").Value = 0.2962807848215612
$ws.Range("M3").Value = 0.029575
$ws.Range("N3").Value = 0.088725
$ws.Range("Q3").Value = 0.0200018091
$ws.Range("R3").Value = 0.1800162819
$ws.Range("S3").Value = 0.2962807848215612
$ws.Range("T3").Value = 0.2962807848215612

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.3513206666666667
$ws.Range("H4").Value = 1.053962
$ws.Range("I4").Value = 0.1539085192605057
$ws.Range("J4").Value = 0.1539085192605057
$ws.Range("M4").Value = 0.029575
$ws.Range("N4").Value = 0.088725
$ws.Range("Q4").Value = 0.01039030871666667
$ws.Range("R4").Value = 0.09351277845
$ws.Range("S4").Value = 0.1539085192605057
$ws.Range("T4").Value = 0.1539085192605057

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.3449053333333333
$ws.Range("H5").Value = 1.034716
$ws.Range("I5").Value = 0.1510980542136751
$ws.Range("J5").Value = 0.1510980542136751
$ws.Range("M5").Value = 0.029575
$ws.Range("N5").Value = 0.088725
$ws.Range("Q5").Value = 0.01020057523333333
$ws.Range("R5").Value = 0.0918051771
$ws.Range("S5").Value = 0.1510980542136751
$ws.Range("T5").Value = 0.1510980542136751

# Row 6
$ws.Range("G6").Value = 0.4544146666666666
$ws.Range("H6").Value = 1.363244
$ws.Range("I6").Value = 0.1990725144082698
$ws.Range("J6").Value = 0.1990725144082698
$ws.Range("M6").Value = 0.029575
$ws.Range("N6").Value = 0.088725
$ws.Range("Q6").Value = 0.01343931376666667
$ws.Range("R6").Value = 0.1209538239
$ws.Range("S6").Value = 0.1990725144082698
$ws.Range("T6").Value = 0.1990725144082698
